$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Re-arrange the Saturday (column G) concert-block content & merges.
#    "Concert call time" moves from G19 -> G20, "Lina Summer Camp..." moves
#    from G20 -> G21, and "After concert refreshment..." moves from G26 -> G27.
# ---------------------------------------------------------------------------
$ws.Range("G15:G18").UnMerge()
$ws.Range("G20:G25").UnMerge()
$ws.Range("G26:G28").UnMerge()

$ws.Range("G19").ClearContents()
$ws.Range("G20").Value = "Concert call time"
$ws.Range("G21").Value = "Lina Summer Camp of Music Students & Friends Concert"
$ws.Range("G26").ClearContents()
$ws.Range("G27").Value = "After concert refreshment `n(Maritime Museum)"

$ws.Range("G15:G19").Merge()
$ws.Range("G21:G26").Merge()
$ws.Range("G27:G28").Merge()

# ---------------------------------------------------------------------------
# 2) Apply consistent borders + centered/wrapped alignment across the whole
#    timetable grid (both the plain body cells and the bold header cells).
# ---------------------------------------------------------------------------
$all = $ws.Range("A1:G30")
$all.Borders.LineStyle = 1
$all.WrapText = $true
$all.HorizontalAlignment = -4108
$all.VerticalAlignment = -4108

# ---------------------------------------------------------------------------
# 3) Restore natural row heights - merging/re-filling wrapped cells can
#    otherwise leave a stale "custom height" on some rows.
# ---------------------------------------------------------------------------
$ws.Rows.Item("1:30").AutoFit()

Write-Host "done"
